# Insert a new data row at row 167 (shifts existing rows 167-252 down to 168-253,
# preserving their content/formatting, and growing the sheet dimension to A1:R253),
# then populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(167).Insert()

$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = 44455
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112023
$ws.Range("G167").Value = "Brócoli"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 800
$ws.Range("K167").Value = 800
$ws.Range("L167").Value = 800
$ws.Range("M167").Value = 800
$ws.Range("N167").Value = "$/unidad"
$ws.Range("O167").Value = "Región Metropolitana"
$ws.Range("P167").Value = 800
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
